$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 51083.668
$ws.Range("I51").Value = 51000.5
$ws.Range("K51").Value = 51000.5
$ws.Range("M51").Value = -50516.5

$ws.Range("H64").Value = 3558.3333
$ws.Range("I64").Value = 3316.6667
$ws.Range("J64").Value = 3800
$ws.Range("K64").Value = 3316.6667
$ws.Range("L64").Value = 3800
$ws.Range("M64").Value = -3068.6667
$ws.Range("N64").Value = -4296

$ws.Range("H67").Value = 3558.3333
$ws.Range("I67").Value = 3316.6667
$ws.Range("J67").Value = 3800
$ws.Range("K67").Value = 3316.6667
$ws.Range("L67").Value = 3800
$ws.Range("M67").Value = -2458.6667
$ws.Range("N67").Value = -5516

$ws.Range("H70").Value = 2062.875
$ws.Range("I70").Value = 2740.6
$ws.Range("J70").Value = 933.3333
$ws.Range("K70").Value = 8221.799999999999
$ws.Range("L70").Value = 2799.9999
$ws.Range("M70").Value = -7951.799999999999
$ws.Range("N70").Value = -3339.9999

$ws.Range("H73").Value = 2062.875
$ws.Range("I73").Value = 2740.6
$ws.Range("J73").Value = 933.3333
$ws.Range("K73").Value = 8221.799999999999
$ws.Range("L73").Value = 2799.9999
$ws.Range("M73").Value = -7285.799999999999
$ws.Range("N73").Value = -4671.9999

$ws.Range("H80").Value = 669.5714
$ws.Range("I80").Value = 520
$ws.Range("J80").Value = 752.6667
$ws.Range("K80").Value = 1560
$ws.Range("L80").Value = 2258.0001
$ws.Range("M80").Value = -562
$ws.Range("N80").Value = -4254.0001

$ws.Range("H83").Value = 669.5714
$ws.Range("I83").Value = 520
$ws.Range("J83").Value = 752.6667
$ws.Range("K83").Value = 4680
$ws.Range("L83").Value = 6774.0003
$ws.Range("M83").Value = 312
$ws.Range("N83").Value = -16758.0003

$ws.Range("H111").Value = 1489.8572
$ws.Range("J111").Value = 1500
$ws.Range("L111").Value = 4500
$ws.Range("N111").Value = -10634

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H68").Value = 16266.667
$ws.Range("I68").Value = 9000
$ws.Range("J68").Value = 19900
$ws.Range("K68").Value = 9000
$ws.Range("L68").Value = 19900
$ws.Range("M68").Value = -8189
$ws.Range("N68").Value = -21522

$ws.Range("H71").Value = 16266.667
$ws.Range("I71").Value = 9000
$ws.Range("J71").Value = 19900
$ws.Range("K71").Value = 27000
$ws.Range("L71").Value = 59700
$ws.Range("M71").Value = -22944
$ws.Range("N71").Value = -67812

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1012904.8
$ws.Range("I86").Value = 1815.1538
$ws.Range("J86").Value = 2327321.5
$ws.Range("K86").Value = 1815.1538
$ws.Range("L86").Value = 2327321.5
$ws.Range("M86").Value = -692.1538
$ws.Range("N86").Value = -2329567.5

$ws.Range("H89").Value = 1012904.8
$ws.Range("I89").Value = 1815.1538
$ws.Range("J89").Value = 2327321.5
$ws.Range("K89").Value = 9075.769
$ws.Range("L89").Value = 11636607.5
$ws.Range("M89").Value = -3459.769
$ws.Range("N89").Value = -11647839.5

$ws.Range("H134").Value = 4837114
$ws.Range("I134").Value = 1883.7222
$ws.Range("J134").Value = 22243942
$ws.Range("K134").Value = 5651.1666
$ws.Range("L134").Value = 66731826
$ws.Range("M134").Value = -3116.1666
$ws.Range("N134").Value = -66736896

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 35698.5
$ws.Range("J70").Value = 35698.5
$ws.Range("L70").Value = 35698.5
$ws.Range("N70").Value = -36328.5

$ws.Range("H73").Value = 35698.5
$ws.Range("J73").Value = 35698.5
$ws.Range("L73").Value = 35698.5
$ws.Range("N73").Value = -37882.5

$ws.Range("H81").Value = 29802.666
$ws.Range("J81").Value = 29802.666
$ws.Range("L81").Value = 29802.666
$ws.Range("N81").Value = -31798.666

$ws.Range("H84").Value = 29802.666
$ws.Range("J84").Value = 29802.666
$ws.Range("L84").Value = 89407.99800000001
$ws.Range("N84").Value = -99391.99800000001

$ws.Range("H132").Value = 22224238
$ws.Range("I132").Value = 1201.7142
$ws.Range("J132").Value = 41669396
$ws.Range("K132").Value = 3605.1426
$ws.Range("L132").Value = 125008188
$ws.Range("M132").Value = -1075.1426
$ws.Range("N132").Value = -125013248

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 41670796
$ws.Range("I5").Value = 53030652
$ws.Range("J5").Value = 18000
$ws.Range("K5").Value = 159091956
$ws.Range("L5").Value = 54000
$ws.Range("M5").Value = -159091844
$ws.Range("N5").Value = -54224

$ws.Range("H86").Value = 200
$ws.Range("I86").Value = 200
$ws.Range("J86").Value = 200
$ws.Range("K86").Value = 600
$ws.Range("L86").Value = 600
$ws.Range("M86").Value = 586
$ws.Range("N86").Value = -2972

$ws.Range("H89").Value = 200
$ws.Range("I89").Value = 200
$ws.Range("J89").Value = 200
$ws.Range("K89").Value = 1800
$ws.Range("L89").Value = 1800
$ws.Range("M89").Value = 4128
$ws.Range("N89").Value = -13656

$ws.Range("H107").Value = 5565527.5
$ws.Range("J107").Value = 6071402.5
$ws.Range("L107").Value = 18214207.5
$ws.Range("N107").Value = -18218047.5

$ws.Range("H131").Value = 867.55554
$ws.Range("I131").Value = 538
$ws.Range("J131").Value = 885.0851
$ws.Range("K131").Value = 1614
$ws.Range("L131").Value = 2655.2553
$ws.Range("M131").Value = 3426
$ws.Range("N131").Value = -12735.2553

$ws.Range("H135").Value = 41670796
$ws.Range("I135").Value = 53030652
$ws.Range("J135").Value = 18000
$ws.Range("K135").Value = 477275868
$ws.Range("L135").Value = 162000
$ws.Range("M135").Value = -477273333
$ws.Range("N135").Value = -167070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7425.5454
$ws.Range("I70").Value = 8721.956
$ws.Range("J70").Value = 4443.8
$ws.Range("K70").Value = 8721.956
$ws.Range("L70").Value = 4443.8
$ws.Range("M70").Value = -8451.956
$ws.Range("N70").Value = -4983.8

$ws.Range("H73").Value = 7425.5454
$ws.Range("I73").Value = 8721.956
$ws.Range("J73").Value = 4443.8
$ws.Range("K73").Value = 8721.956
$ws.Range("L73").Value = 4443.8
$ws.Range("M73").Value = -7785.956
$ws.Range("N73").Value = -6315.8

$ws.Range("H113").Value = 3400
$ws.Range("I113").Value = 3400
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3400
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1230
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 62502350
$ws.Range("I61").Value = 2336
$ws.Range("K61").Value = 2336
$ws.Range("M61").Value = -2134

$ws.Range("H113").Value = 62502350
$ws.Range("I113").Value = 2336
$ws.Range("K113").Value = 2336
$ws.Range("M113").Value = -166

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 34292.75
$ws.Range("J68").Value = 34292.75
$ws.Range("L68").Value = 34292.75
$ws.Range("N68").Value = -35914.75

$ws.Range("H71").Value = 34292.75
$ws.Range("J71").Value = 34292.75
$ws.Range("L71").Value = 102878.25
$ws.Range("N71").Value = -110990.25

$ws.Range("H75").Value = 18912.8
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 18912.8
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 18912.8
$ws.Range("N75").Value = -20784.8
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 18912.8
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 18912.8
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 56738.39999999999
$ws.Range("N78").Value = -66098.39999999999
$ws.Range("M78").ClearContents()

$ws.Range("H99").Value = 25000
$ws.Range("J99").Value = 25000
$ws.Range("L99").Value = 25000
$ws.Range("N99").Value = -30990

$ws.Range("H136").Value = 1469.8363
$ws.Range("I136").Value = 935.9286
$ws.Range("K136").Value = 2807.7858
$ws.Range("M136").Value = -257.7857999999997

